# Natmi following Dr Hou advice
# Update St6gal1-Cd22 LR-pair results: recomputed stats for existing
# Sending-cluster x M2 rows, and add matching Sending-cluster x ECs rows
# for each of ECs, FAPs, M2, sCs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$row2 = @("ECs", "St6gal1", "Cd22", "ECs", 3, 1, 3.343101333333333, 10.029304, 0.2255238790820294, 0.2255238790820294, 1, 0.3333333333333333, 0.5751623333333333, 1.725487, 0.08167574423707133, 0.08167574423707133, 1.922825963449778, 17.305433671048, 0.01841983066725603, 0.01841983066725603)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

# Row 3: ECs -> M2
$row3 = @("ECs", "St6gal1", "Cd22", "M2", 3, 1, 3.343101333333333, 10.029304, 0.2255238790820294, 0.2255238790820294, 3, 1, 6.466858999999999, 19.400577, 0.9183242557629286, 0.9183242557629286, 21.61936494537866, 194.574284508408, 0.2071040484147733, 0.2071040484147734)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}

# Row 4: FAPs -> ECs
$row4 = @("FAPs", "St6gal1", "Cd22", "ECs", 3, 1, 3.027316, 9.081948, 0.2042211645575085, 0.2042211645575085, 1, 0.3333333333333333, 0.5751623333333333, 1.725487, 0.08167574423707133, 0.08167574423707133, 1.741198134297333, 15.670783208676, 0.01667991560419593, 0.01667991560419593)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $row4[$i]
}

# Row 5: FAPs -> M2
$row5 = @("FAPs", "St6gal1", "Cd22", "M2", 3, 1, 3.027316, 9.081948, 0.2042211645575085, 0.2042211645575085, 3, 1, 6.466858999999999, 19.400577, 0.9183242557629286, 0.9183242557629286, 19.577225720444, 176.195031483996, 0.1875412489533126, 0.1875412489533126)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = $row5[$i]
}

# Row 6: M2 -> ECs
$row6 = @("M2", "St6gal1", "Cd22", "ECs", 3, 1, 5.111358333333333, 15.334075, 0.3448095776271982, 0.3448095776271983, 1, 0.3333333333333333, 0.5751623333333333, 1.725487, 0.08167574423707133, 0.08167574423707133, 2.939860785502778, 26.458747069525, 0.02816257887277163, 0.02816257887277164)
for ($i = 0; $i -lt $row6.Length; $i++) {
    $ws.Cells.Item(6, $i + 1).Value = $row6[$i]
}

# Row 7: M2 -> M2
$row7 = @("M2", "St6gal1", "Cd22", "M2", 3, 1, 5.111358333333333, 15.334075, 0.3448095776271982, 0.3448095776271983, 3, 1, 6.466858999999999, 19.400577, 0.9183242557629286, 0.9183242557629286, 33.05443364014167, 297.489902761275, 0.3166469987544266, 0.3166469987544266)
for ($i = 0; $i -lt $row7.Length; $i++) {
    $ws.Cells.Item(7, $i + 1).Value = $row7[$i]
}

# Row 8: sCs -> ECs
$row8 = @("sCs", "St6gal1", "Cd22", "ECs", 3, 1, 3.341937666666667, 10.025813, 0.2254453787332639, 0.2254453787332639, 1, 0.3333333333333333, 0.5751623333333333, 1.725487, 0.08167574423707133, 0.08167574423707133, 1.922156666214556, 17.299409995931, 0.01841341909284774, 0.01841341909284774)
for ($i = 0; $i -lt $row8.Length; $i++) {
    $ws.Cells.Item(8, $i + 1).Value = $row8[$i]
}

# Row 9: sCs -> M2
$row9 = @("sCs", "St6gal1", "Cd22", "M2", 3, 1, 3.341937666666667, 10.025813, 0.2254453787332639, 0.2254453787332639, 3, 1, 6.466858999999999, 19.400577, 0.9183242557629286, 0.9183242557629286, 21.61183967712233, 194.506557094101, 0.2070319596404161, 0.2070319596404162)
for ($i = 0; $i -lt $row9.Length; $i++) {
    $ws.Cells.Item(9, $i + 1).Value = $row9[$i]
}

